# Updated symbol list on Sat Dec 31 11:42:23 UTC 2022 with GitHub Actions
# All Price (column D) values are stored as literal text in the workbook,
# so numeric-looking values are written with a leading apostrophe to force
# Excel to keep them as text (preserving exact formatting / trailing zeros)
# instead of coercing them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'26.14"

$ws.Range("D4").Value = "'5.092"

$ws.Range("D5").Value = "'0.05595"

$ws.Range("D6").Value = "'6.497"

$ws.Range("D8").Value = "'0.8109"

$ws.Range("D9").Value = "'0.8424"

$ws.Range("D10").Value = "'0.1345"

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = "'0.03229"
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = "'0.02770"
$ws.Range("E12").Value = '11BitrueCoinBTR'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = "'0.09401"
$ws.Range("E13").Value = '12BitMartTokenBMX'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = "'0.001512"
$ws.Range("E14").Value = '13BitForexTokenBF'

$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = "'0.0005990"
$ws.Range("E15").Value = '14OneONE'

$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = "'0.006115"
$ws.Range("E16").Value = '15TigerCashTCH'

$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = "'3.558"
$ws.Range("E17").Value = '16LEOLEO'

$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = "'2.118"
$ws.Range("E18").Value = '17BTSETokenBTSE'

$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = "'0.3183"
$ws.Range("E19").Value = '18BitpandaEcosystemTokenBEST'

$ws.Range("B20").Value = 'MandalaExchangeToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D20").Value = "'0.07005"
$ws.Range("E20").Value = '19MandalaExchangeTokenMDX'

$ws.Range("D22").Value = "'3.753"

$ws.Range("D23").Value = "'0.04713"

$ws.Range("D26").Value = "'0.004611"

$ws.Range("D40").Value = "'0.03654"

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = "'0.1358"
$ws.Range("E41").Value = '40BKEXTokenBKKBestin24h'

$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.002660"
$ws.Range("E42").Value = '41CEJICEJI'

$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = "'0.006120"
$ws.Range("E43").Value = '42KickTokenKICK'

$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'

$ws.Range("D48").Value = "'0.002053"
